$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(16, 8).Value = 36505  # H16
$ws.Cells.Item(16, 9).Value = 0  # I16
$ws.Cells.Item(16, 10).Value = 36505  # J16
$ws.Cells.Item(16, 11).Value = 0  # K16
$ws.Cells.Item(16, 12).Value = 36505  # L16
$ws.Cells.Item(16, 13).ClearContents()  # M16
$ws.Cells.Item(16, 14).Value = -36965  # N16
$ws.Cells.Item(21, 8).Value = 9572.333000000001  # H21
$ws.Cells.Item(21, 9).Value = 9572.333000000001  # I21
$ws.Cells.Item(21, 11).Value = 9572.333000000001  # K21
$ws.Cells.Item(21, 13).Value = -9104.333000000001  # M21
$ws.Cells.Item(23, 8).Value = 9572.333000000001  # H23
$ws.Cells.Item(23, 9).Value = 9572.333000000001  # I23
$ws.Cells.Item(23, 11).Value = 9572.333000000001  # K23
$ws.Cells.Item(23, 13).Value = -9338.333000000001  # M23
$ws.Cells.Item(34, 8).Value = 24272  # H34
$ws.Cells.Item(34, 9).Value = 1803.1428  # I34
$ws.Cells.Item(34, 10).Value = 76699.336  # J34
$ws.Cells.Item(34, 11).Value = 1803.1428  # K34
$ws.Cells.Item(34, 12).Value = 76699.336  # L34
$ws.Cells.Item(34, 13).Value = -1600.1428  # M34
$ws.Cells.Item(34, 14).Value = -77105.336  # N34
$ws.Cells.Item(36, 8).Value = 24272  # H36
$ws.Cells.Item(36, 9).Value = 1803.1428  # I36
$ws.Cells.Item(36, 10).Value = 76699.336  # J36
$ws.Cells.Item(36, 11).Value = 1803.1428  # K36
$ws.Cells.Item(36, 12).Value = 76699.336  # L36
$ws.Cells.Item(36, 13).Value = -1088.1428  # M36
$ws.Cells.Item(36, 14).Value = -78129.336  # N36
$ws.Cells.Item(47, 8).Value = 80049.336  # H47
$ws.Cells.Item(47, 10).Value = 80049.336  # J47
$ws.Cells.Item(47, 12).Value = 80049.336  # L47
$ws.Cells.Item(47, 14).Value = -81993.336  # N47
$ws.Cells.Item(107, 8).Value = 407.38095  # H107
$ws.Cells.Item(107, 9).Value = 365.21054  # I107
$ws.Cells.Item(107, 10).Value = 808  # J107
$ws.Cells.Item(107, 11).Value = 365.21054  # K107
$ws.Cells.Item(107, 12).Value = 808  # L107
$ws.Cells.Item(107, 13).Value = 1554.78946  # M107
$ws.Cells.Item(107, 14).Value = -4648  # N107
$ws.Cells.Item(125, 8).Value = 1365.1111  # H125
$ws.Cells.Item(125, 9).Value = 675  # I125
$ws.Cells.Item(125, 10).Value = 1562.2858  # J125
$ws.Cells.Item(125, 11).Value = 6075  # K125
$ws.Cells.Item(125, 12).Value = 14060.5722  # L125
$ws.Cells.Item(125, 13).Value = -3615  # M125
$ws.Cells.Item(125, 14).Value = -18980.5722  # N125
$ws.Cells.Item(132, 8).Value = 6494958.5  # H132
$ws.Cells.Item(132, 9).Value = 8404670  # I132
$ws.Cells.Item(132, 11).Value = 25214010  # K132
$ws.Cells.Item(132, 13).Value = -25211480  # M132
$ws.Cells.Item(138, 8).Value = 1385.1791  # H138
$ws.Cells.Item(138, 9).Value = 817.2895  # I138
$ws.Cells.Item(138, 10).Value = 2129.3103  # J138
$ws.Cells.Item(138, 11).Value = 2451.8685  # K138
$ws.Cells.Item(138, 12).Value = 6387.9309  # L138
$ws.Cells.Item(138, 13).Value = 2688.1315  # M138
$ws.Cells.Item(138, 14).Value = -16667.9309  # N138

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 964.8  # H2
$ws.Cells.Item(2, 9).Value = 635.1111  # I2
$ws.Cells.Item(2, 10).Value = 1459.3334  # J2
$ws.Cells.Item(2, 11).Value = 635.1111  # K2
$ws.Cells.Item(2, 12).Value = 1459.3334  # L2
$ws.Cells.Item(2, 13).Value = -522.1111  # M2
$ws.Cells.Item(2, 14).Value = -1685.3334  # N2
$ws.Cells.Item(45, 8).Value = 970.05  # H45
$ws.Cells.Item(45, 9).Value = 1036.5  # I45
$ws.Cells.Item(45, 11).Value = 1036.5  # K45
$ws.Cells.Item(45, 13).Value = -659.5  # M45
$ws.Cells.Item(61, 8).Value = 2514.9524  # H61
$ws.Cells.Item(61, 9).Value = 2074  # I61
$ws.Cells.Item(61, 10).Value = 3000  # J61
$ws.Cells.Item(61, 11).Value = 2074  # K61
$ws.Cells.Item(61, 12).Value = 3000  # L61
$ws.Cells.Item(61, 13).Value = -1862  # M61
$ws.Cells.Item(61, 14).Value = -3424  # N61
$ws.Cells.Item(74, 8).Value = 1906.0714  # H74
$ws.Cells.Item(74, 9).Value = 1932.0834  # I74
$ws.Cells.Item(74, 11).Value = 1932.0834  # K74
$ws.Cells.Item(74, 13).Value = -1058.0834  # M74
$ws.Cells.Item(77, 8).Value = 1906.0714  # H77
$ws.Cells.Item(77, 9).Value = 1932.0834  # I77
$ws.Cells.Item(77, 11).Value = 9660.416999999999  # K77
$ws.Cells.Item(77, 13).Value = -5292.416999999999  # M77
$ws.Cells.Item(116, 8).Value = 964.8  # H116
$ws.Cells.Item(116, 9).Value = 635.1111  # I116
$ws.Cells.Item(116, 10).Value = 1459.3334  # J116
$ws.Cells.Item(116, 11).Value = 635.1111  # K116
$ws.Cells.Item(116, 12).Value = 1459.3334  # L116
$ws.Cells.Item(116, 13).Value = 1658.8889  # M116
$ws.Cells.Item(116, 14).Value = -6047.3334  # N116
$ws.Cells.Item(136, 8).Value = 2514.9524  # H136
$ws.Cells.Item(136, 9).Value = 2074  # I136
$ws.Cells.Item(136, 10).Value = 3000  # J136
$ws.Cells.Item(136, 11).Value = 6222  # K136
$ws.Cells.Item(136, 12).Value = 9000  # L136
$ws.Cells.Item(136, 13).Value = -3672  # M136
$ws.Cells.Item(136, 14).Value = -14100  # N136

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 964.8  # H3
$ws.Cells.Item(3, 9).Value = 635.1111  # I3
$ws.Cells.Item(3, 10).Value = 1459.3334  # J3
$ws.Cells.Item(3, 11).Value = 635.1111  # K3
$ws.Cells.Item(3, 12).Value = 1459.3334  # L3
$ws.Cells.Item(3, 13).Value = -521.1111  # M3
$ws.Cells.Item(3, 14).Value = -1687.3334  # N3
$ws.Cells.Item(20, 8).Value = 1930.4286  # H20
$ws.Cells.Item(20, 9).Value = 1903.3914  # I20
$ws.Cells.Item(20, 10).Value = 1982.25  # J20
$ws.Cells.Item(20, 11).Value = 1903.3914  # K20
$ws.Cells.Item(20, 12).Value = 1982.25  # L20
$ws.Cells.Item(20, 13).Value = -1656.3914  # M20
$ws.Cells.Item(20, 14).Value = -2476.25  # N20
$ws.Cells.Item(52, 8).Value = 30390  # H52
$ws.Cells.Item(52, 10).Value = 30390  # J52
$ws.Cells.Item(52, 12).Value = 30390  # L52
$ws.Cells.Item(52, 14).Value = -30916  # N52
$ws.Cells.Item(107, 8).Value = 821.1  # H107
$ws.Cells.Item(107, 9).Value = 821.1  # I107
$ws.Cells.Item(107, 11).Value = 821.1  # K107
$ws.Cells.Item(107, 13).Value = 1098.9  # M107
$ws.Cells.Item(118, 8).Value = 32000  # H118
$ws.Cells.Item(118, 10).Value = 32000  # J118
$ws.Cells.Item(118, 12).Value = 32000  # L118
$ws.Cells.Item(118, 14).Value = -35314  # N118
$ws.Cells.Item(119, 8).Value = 29000  # H119
$ws.Cells.Item(119, 10).Value = 29000  # J119
$ws.Cells.Item(119, 12).Value = 29000  # L119
$ws.Cells.Item(119, 14).Value = -38676  # N119
$ws.Cells.Item(120, 8).Value = 30000  # H120
$ws.Cells.Item(120, 10).Value = 30000  # J120
$ws.Cells.Item(120, 12).Value = 30000  # L120
$ws.Cells.Item(120, 14).Value = -39676  # N120
$ws.Cells.Item(121, 8).Value = 30390  # H121
$ws.Cells.Item(121, 10).Value = 30390  # J121
$ws.Cells.Item(121, 12).Value = 30390  # L121
$ws.Cells.Item(121, 14).Value = -33884  # N121

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1390  # H16
$ws.Cells.Item(16, 9).Value = 1072.5  # I16
$ws.Cells.Item(16, 10).Value = 1571.4286  # J16
$ws.Cells.Item(16, 11).Value = 1072.5  # K16
$ws.Cells.Item(16, 12).Value = 1571.4286  # L16
$ws.Cells.Item(16, 13).Value = -785.5  # M16
$ws.Cells.Item(16, 14).Value = -2145.4286  # N16
$ws.Cells.Item(31, 8).Value = 4002284.2  # H31
$ws.Cells.Item(31, 9).Value = 2600.5186  # I31
$ws.Cells.Item(31, 10).Value = 8697565  # J31
$ws.Cells.Item(31, 11).Value = 2600.5186  # K31
$ws.Cells.Item(31, 12).Value = 8697565  # L31
$ws.Cells.Item(31, 13).Value = -2305.5186  # M31
$ws.Cells.Item(31, 14).Value = -8698155  # N31
$ws.Cells.Item(34, 8).Value = 4002284.2  # H34
$ws.Cells.Item(34, 9).Value = 2600.5186  # I34
$ws.Cells.Item(34, 10).Value = 8697565  # J34
$ws.Cells.Item(34, 11).Value = 2600.5186  # K34
$ws.Cells.Item(34, 12).Value = 8697565  # L34
$ws.Cells.Item(34, 13).Value = -2398.5186  # M34
$ws.Cells.Item(34, 14).Value = -8697969  # N34
$ws.Cells.Item(86, 8).Value = 5000  # H86
$ws.Cells.Item(86, 9).Value = 0  # I86
$ws.Cells.Item(86, 11).Value = 0  # K86
$ws.Cells.Item(86, 13).ClearContents()  # M86
$ws.Cells.Item(89, 8).Value = 5000  # H89
$ws.Cells.Item(89, 9).Value = 0  # I89
$ws.Cells.Item(89, 11).Value = 0  # K89
$ws.Cells.Item(89, 13).ClearContents()  # M89
$ws.Cells.Item(113, 8).Value = 1390  # H113
$ws.Cells.Item(113, 9).Value = 1072.5  # I113
$ws.Cells.Item(113, 10).Value = 1571.4286  # J113
$ws.Cells.Item(113, 11).Value = 1072.5  # K113
$ws.Cells.Item(113, 12).Value = 1571.4286  # L113
$ws.Cells.Item(113, 13).Value = 1097.5  # M113
$ws.Cells.Item(113, 14).Value = -5911.4286  # N113

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 638.4783  # H107
$ws.Cells.Item(107, 9).Value = 599.1177  # I107
$ws.Cells.Item(107, 10).Value = 750  # J107
$ws.Cells.Item(107, 11).Value = 599.1177  # K107
$ws.Cells.Item(107, 12).Value = 750  # L107
$ws.Cells.Item(107, 13).Value = 1320.8823  # M107
$ws.Cells.Item(107, 14).Value = -4590  # N107

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 4679.4688  # H136
$ws.Cells.Item(136, 9).Value = 6311.737  # I136
$ws.Cells.Item(136, 10).Value = 2293.8462  # J136
$ws.Cells.Item(136, 11).Value = 18935.211  # K136
$ws.Cells.Item(136, 12).Value = 6881.5386  # L136
$ws.Cells.Item(136, 13).Value = -16385.211  # M136
$ws.Cells.Item(136, 14).Value = -11981.5386  # N136

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 461.05264  # H100
$ws.Cells.Item(100, 9).Value = 436.66666  # I100
$ws.Cells.Item(100, 10).Value = 502.85715  # J100
$ws.Cells.Item(100, 11).Value = 873.33332  # K100
$ws.Cells.Item(100, 12).Value = 1005.7143  # L100
$ws.Cells.Item(100, 13).Value = -332.33332  # M100
$ws.Cells.Item(100, 14).Value = -2087.7143  # N100
$ws.Cells.Item(113, 8).Value = 775  # H113
$ws.Cells.Item(113, 9).Value = 781.8182  # I113
$ws.Cells.Item(113, 11).Value = 2345.4546  # K113
$ws.Cells.Item(113, 13).Value = -175.4546  # M113
$ws.Cells.Item(132, 8).Value = 1522.6786  # H132
$ws.Cells.Item(132, 9).Value = 1349.1305  # I132
$ws.Cells.Item(132, 10).Value = 2321  # J132
$ws.Cells.Item(132, 11).Value = 4047.3915  # K132
$ws.Cells.Item(132, 12).Value = 6963  # L132
$ws.Cells.Item(132, 13).Value = -1517.3915  # M132
$ws.Cells.Item(132, 14).Value = -12023  # N132
$ws.Cells.Item(136, 8).Value = 5340.077  # H136
$ws.Cells.Item(136, 9).Value = 7308  # I136
$ws.Cells.Item(136, 10).Value = 912.25  # J136
$ws.Cells.Item(136, 11).Value = 21924  # K136
$ws.Cells.Item(136, 12).Value = 2736.75  # L136
$ws.Cells.Item(136, 13).Value = -19374  # M136
